# generator: const range list added
#
# Adds a "SAND_PROP_RANGE" (min/max) example to the "test" sheet (columns
# G/H/I), and fills in the actual min/max range values for the "time" sheet
# rows (month/day/hour/minute/sec/msec/utc). Finishes with the "test" sheet
# selected/active (matching the saved workbook view).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "time" sheet (4th tab): fill in H/I (min/max) range columns.
# ---------------------------------------------------------------------
$wsTime = $wb.Worksheets.Item(4)

$wsTime.Range("H4").Value = 1      # month min
$wsTime.Range("I4").Value = 12     # month max

$wsTime.Range("H5").Value = 1      # day min
$wsTime.Range("I5").Value = 31     # day max

$wsTime.Range("H6").Value = 0      # hour min
$wsTime.Range("I6").Value = 23     # hour max

$wsTime.Range("H7").Value = 0      # minute min
$wsTime.Range("I7").Value = 59     # minute max

$wsTime.Range("H8").Value = 0      # sec min
$wsTime.Range("I8").Value = 59     # sec max

$wsTime.Range("H9").Value = 0      # msec min
$wsTime.Range("I9").Value = 999    # msec max

$wsTime.Range("H11").Value = -12   # utc min
$wsTime.Range("I11").Value = 14    # utc max

$wsTime.Range("I12").Select()

# ---------------------------------------------------------------------
# 2) "test" sheet (10th/last tab): const range / literal examples.
# ---------------------------------------------------------------------
$wsTest = $wb.Worksheets.Item(10)

# Widen column G to match column F so the new example values are legible.
$wsTest.Columns.Item(7).ColumnWidth = $wsTest.Columns.Item(6).ColumnWidth

# Values are entered in this order so the shared-string table comes out in
# the same order as the original edit.
$wsTest.Range("G3").Value = "{17,0,4,2,5}"        # arr_u8  sample literal
$wsTest.Range("G13").Value = "`"Hello world!`""   # arr_char sample literal
$wsTest.Range("G11").Value = "{15.4f, 71.524f}"   # arr_float sample literal
$wsTest.Range("H11").Value = "'-12.5f"            # arr_float min (quote-prefixed text)
$wsTest.Range("I11").Value = "255.44f"            # arr_float max
$wsTest.Range("G12").Value = "{15.4f,71.524f}"    # arr_double sample literal

$wsTest.Range("H13").Select()

# Make "test" the active/selected sheet (matches the saved workbook view).
$wsTest.Activate()
